$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K4").Value = 19
$ws.Range("J5").Value = 1.14
$ws.Range("K5").Value = 5.5
$ws.Range("G6").Value = 3.5
$ws.Range("I6").Value = 2.3
$ws.Range("T6").Value = 8
$ws.Range("U6").Value = 15
$ws.Range("AA6").Value = 5.5
$ws.Range("AF6").Value = 10
$ws.Range("G8").Value = 1.38
$ws.Range("H8").Value = 4.75
$ws.Range("I8").Value = 8
$ws.Range("T8").Value = 8.5
$ws.Range("U8").Value = 7.5
$ws.Range("W8").Value = 9.5
$ws.Range("AC8").Value = 51
$ws.Range("AE8").Value = 23
$ws.Range("AG8").Value = 23
$ws.Range("AI8").Value = 51
$ws.Range("AJ8").Value = 51
$ws.Range("G11").Value = 1.95
$ws.Range("I11").Value = 3.7
$ws.Range("N11").Value = 2.15
$ws.Range("O11").Value = 1.55
$ws.Range("V11").Value = 8.75
$ws.Range("W11").Value = 16.5
$ws.Range("X11").Value = 17.5
$ws.Range("Z11").Value = 7.7
$ws.Range("AA11").Value = 6.3
$ws.Range("AE11").Value = 8.75
$ws.Range("AF11").Value = 18.5
$ws.Range("AI11").Value = 40
$ws.Range("AJ11").Value = 55
$ws.Range("H12").Value = 3.15
$ws.Range("X12").Value = 18
$ws.Range("AE12").Value = 9.25
$ws.Range("AH12").Value = 90
$ws.Range("G13").Value = 2.1
$ws.Range("H13").Value = 3.2
$ws.Range("I13").Value = 3.6
$ws.Range("L13").Value = 1.44
$ws.Range("M13").Value = 2.63
$ws.Range("N13").Value = 2.35
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.5
$ws.Range("R13").Value = 2.05
$ws.Range("S13").Value = 1.7
$ws.Range("T13").Value = 6
$ws.Range("U13").Value = 9
$ws.Range("X13").Value = 21
$ws.Range("Z13").Value = 7
$ws.Range("AB13").Value = 19
$ws.Range("AC13").Value = 67
$ws.Range("AD13").Value = 501
$ws.Range("G17").Value = 2.2
$ws.Range("I17").Value = 3.6
$ws.Range("K17").Value = 7.5
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 10
$ws.Range("W17").Value = 21
$ws.Range("AI17").Value = 29
$ws.Range("N19").Value = 1.73
$ws.Range("O19").Value = 2.08
$ws.Range("N23").Value = 2
$ws.Range("O23").Value = 1.8
$ws.Range("N24").Value = 1.85
$ws.Range("O24").Value = 1.95
$ws.Range("G26").Value = 3.6
$ws.Range("I26").Value = 1.8
$ws.Range("J26").Value = 1.01
$ws.Range("K26").Value = 17
$ws.Range("AA26").Value = 8
$ws.Range("AE26").Value = 10
$ws.Range("J27").Value = 1.01
$ws.Range("K27").Value = 15
$ws.Range("G30").Value = 1.83
$ws.Range("I30").Value = 3.9
$ws.Range("W30").Value = 17
$ws.Range("AF30").Value = 21
$ws.Range("N33").Value = 1.93
$ws.Range("O33").Value = 1.93
$ws.Range("N34").Value = 1.62
$ws.Range("O34").Value = 2.25
